# MHD2-137: Support for creating TWIST Haem assay reports
#
# 1) "RESULTS SUMMARY_IN" -> "RESULTS_SUMMARY_IN" (space becomes underscore)
# 2) Tidy up the "Test Methodology" paragraph: the text itself is unchanged,
#    but re-typing each sentence over the existing (spell-check-marked) runs
#    clears the stale w:proofErr wavy-underline markers left over from the
#    older wording.
# 3) Report date bumped from 16-Sep-2024 to 17-Sep-2024.

$d = $word.ActiveDocument
$wdReplaceOne = 1
$wdFindContinue = 1

# --- 1. RESULTS SUMMARY_IN -> RESULTS_SUMMARY_IN -------------------------
$rng = $d.Content
[void]$rng.Find.Execute(
    "RESULTS SUMMARY_IN", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "RESULTS_SUMMARY_IN", $wdReplaceOne)

# --- 2. Test methodology paragraph clean-up -------------------------------
$part1 = "panel (Peter MacCallum Cancer Centre AllHaem DNA Twist v1, design ID TE-98899881) and sequenced on an Illumina NovaSeq 6000 with 150 bp paired end reads."
$rng = $d.Content
[void]$rng.Find.Execute(
    $part1, $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, $part1, $wdReplaceOne)

$part2 = "A custom Seqliner/Nextflow-based analysis pipeline is used to generate aligned reads and call variants (single nucleotide variants and short insertions or deletions) against the hg19 human reference genome."
$rng = $d.Content
[void]$rng.Find.Execute(
    $part2, $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, $part2, $wdReplaceOne)

$part3 = "Variants are analysed using PathOS software (Peter Mac)"
$rng = $d.Content
[void]$rng.Find.Execute(
    $part3, $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, $part3, $wdReplaceOne)

# --- 3. Report date ---------------------------------------------------------
$rng = $d.Content
[void]$rng.Find.Execute(
    "16-Sep-2024", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "17-Sep-2024", $wdReplaceOne)
